$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove obsolete rows (old users 1234567890123450, 3201234567890000,
# 12345678900000, 2345678990112444, 1245367800112234, 9801234567819235),
# keeping only the header row and the row for "Nia" (originally row 5,
# which becomes the new row 2). Delete the rows below "Nia" first so the
# row numbers above are not affected by the shift.
$ws.Rows("6:8").Delete()
$ws.Rows("2:4").Delete()

# Add the new "Password" column header, copying the bold/border/centered
# style used by the other header cells (A1:C1).
$ws.Range("A1").Copy($ws.Range("D1"))
$ws.Range("D1").Value = "Password"

# Update the remaining data row (previously row 5 / "Nia") with her
# fuller name (NIK in A2 and the blank Plat in B2 already carry over
# correctly from the old row), and add the new account password.
$ws.Range("C2").Value = "Nia Rahmadani"
$ws.Range("D2").Value = "Akun_nia21"
